# RPA datasets push 2024-07-17
# Sheet 1 (01_IB전략컨설팅부): remove the two rows for "KB제28호스팩" and
# "아이씨티케이" (old rows 20 and 21). This shifts the former row 22
# ("신한글로벌액티브리츠") up to become the new row 20.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(20).Delete()
$ws1.Rows.Item(20).Delete()

# Sheet 2 (02_38커뮤니케이션(최근일자기준)): insert a new row for "엠83"
# right after the existing "에이치이엠파마(구.에이치이엠)" row, and drop the
# last row ("하스") so the sheet keeps the same number of rows.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(3).Insert()
$ws2.Range("A3").Value = "엠83"
$ws2.Range("B3").Value = "2024.08.01~08.07"
$ws2.Range("C3").Value = "11,000~13,000"
$ws2.Range("D3").Value = "-"
$ws2.Range("E3").Value = 16500
$ws2.Range("F3").Value = "신영증권,유진투자증권"

$ws2.Rows.Item(22).Delete()
